$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.297.72'
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").Value = '1.680.66'
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.50'
$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5277'
$ws.Range("E6").Value = '  +2.99%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("E8").Value = '  +2.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06492'
$ws.Range("E9").Value = '  +1.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.00'
$ws.Range("E10").Value = '  +1.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07534'
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("D12").Value = '1.687.98'
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.535'
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5812'
$ws.Range("E14").Value = '  -0.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008514'
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.68'
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").Value = '26.334.84'
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.925'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.89'
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '190.49'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.202'
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '145.46'
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.806'
$ws.Range("E25").Value = '  +2.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1245'
$ws.Range("E26").Value = '  +4.09%  '
$ws.Range("E27").Value = '  +1.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06526'
$ws.Range("E28").Value = '  +2.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.356'
$ws.Range("E29").Value = '  +4.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.332'
$ws.Range("E30").Value = '  +0.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.603'
$ws.Range("E31").Value = '  +2.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.596'
$ws.Range("E32").Value = '  +2.10%  '
$ws.Range("E33").Value = '  +1.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.034'
$ws.Range("E34").Value = '  +1.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6242'
$ws.Range("E35").Value = '  +2.50%  '
$ws.Range("E36").Value = '  +1.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.733'
$ws.Range("E37").Value = '  +2.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.449'
$ws.Range("E38").Value = '  +4.63%  '
$ws.Range("D39").Value = '1.113.02'
$ws.Range("E39").Value = '  +2.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01626'
$ws.Range("E40").Value = '  +1.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8765'
$ws.Range("E41").Value = '  +1.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.015'
$ws.Range("E42").Value = '  +0.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.78'
$ws.Range("E43").Value = '  -0.63%  '
$ws.Range("D44").Value = '1.830.25'
$ws.Range("E44").Value = '  +0.77%  '
$ws.Range("E45").Value = '  -0.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.08'
$ws.Range("E46").Value = '  +1.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.173'
$ws.Range("E47").Value = '  +0.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.008'
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05273'
$ws.Range("E49").Value = '  +1.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.091'
$ws.Range("E50").Value = '  +3.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4291'
$ws.Range("E51").Value = '  -0.07%  '
